$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "UpcomingRequirement Test: STARTED"
$ws.Range("A2").Value = "UpcomingRequirement Test Case: startBrowser Test Method: SUCCESS"
$ws.Range("A3").Value = "UpcomingRequirement Test Case: loadUpcomingBikesPage Test Method: SUCCESS"
$ws.Range("A4").Value = "UpcomingRequirement Test Case: navigateToUpcomingBikesPage Test Method: SUCCESS"
$ws.Range("A5").Value = "UpcomingRequirement Test Case: outputDisplay Test Method: SUCCESS"
$ws.Range("A6").Value = "UpcomingRequirement Test: ENDED"
